# Strip the etapa-specific suffix (_ELE, _ENI, _EOC, _ENIA) from the
# "responsible" column (C) values, e.g. "CamiloR_ELE" -> "CamiloR".
# This consolidates duplicate names across etapas so the dashboard /
# charts group by the bare responsible name instead of per-etapa variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -ne $null) {
        $updated = $current -replace '_(ELE|ENI|EOC|ENIA)$', ''
        if ($updated -ne $current) {
            $cell.Value2 = $updated
        }
    }
}
